$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173, shifting existing rows 173-228 down to 174-229
$ws.Rows.Item(173).Insert()

# Populate the new row 173 with values
$ws.Cells.Item(173, 1).Value = 3
$ws.Cells.Item(173, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44524
$ws.Cells.Item(173, 5).Value = 5
$ws.Cells.Item(173, 6).Value = 100112012
$ws.Cells.Item(173, 7).Value = "Espinaca"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 120
$ws.Cells.Item(173, 11).Value = 3000
$ws.Cells.Item(173, 12).Value = 3000
$ws.Cells.Item(173, 13).Value = 3000
$ws.Cells.Item(173, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(173, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(173, 16).Value = 1000
$ws.Cells.Item(173, 17).Value = 3
$ws.Cells.Item(173, 18).Value = "Hortaliza"
